# get_dash/budget.xlsx bug fixes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet" to "Sheet1"
$ws.Name = "Sheet1"

# Clear out the old table (A1:D4) entirely, including the stray empty row 4
$ws.Range("A1:D4").Clear()

# Write the new header row
$ws.Range("B1").Value = "budget"
$ws.Range("C1").Value = "userId"

# Write the new data rows
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 3000
$ws.Range("C2").Value = "Dixon3220"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 400
$ws.Range("C3").Value = "ntu003"

# Apply bold + thin box border + center/top alignment to the header cells and
# the A-column flag cells, matching the new styling introduced in this revision.
$styledRanges = @("B1:C1", "A2:A3")
foreach ($addr in $styledRanges) { $ws.Range($addr).Font.Bold = $true }
foreach ($addr in $styledRanges) { $ws.Range($addr).Borders.LineStyle = 1 }
foreach ($addr in $styledRanges) { $ws.Range($addr).HorizontalAlignment = -4108 }
foreach ($addr in $styledRanges) { $ws.Range($addr).VerticalAlignment = -4160 }

Write-Host "edit complete"
